$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "86-31=55"
$t.Cell(1,2).Range.Text = "75-13=62"
$t.Cell(1,3).Range.Text = "36+4=40"
$t.Cell(1,4).Range.Text = "23+28=51"
$t.Cell(1,5).Range.Text = "43+43=86"

$t.Cell(2,1).Range.Text = "65-42=23"
$t.Cell(2,2).Range.Text = "57-16=41"
$t.Cell(2,3).Range.Text = "69-48=21"
$t.Cell(2,4).Range.Text = "56-4=52"
$t.Cell(2,5).Range.Text = "5+5=10"

$t.Cell(3,1).Range.Text = "33-24=9"
$t.Cell(3,2).Range.Text = "16+6=22"
$t.Cell(3,3).Range.Text = "80-22=58"
$t.Cell(3,4).Range.Text = "93-53=40"
$t.Cell(3,5).Range.Text = "16+31=47"

$t.Cell(4,1).Range.Text = "70-53=17"
$t.Cell(4,2).Range.Text = "95-93=2"
$t.Cell(4,3).Range.Text = "50+31=81"
$t.Cell(4,4).Range.Text = "64-56=8"
$t.Cell(4,5).Range.Text = "79-3=76"

$t.Cell(5,1).Range.Text = "96+0=96"
$t.Cell(5,2).Range.Text = "80-36=44"
$t.Cell(5,3).Range.Text = "59-57=2"
$t.Cell(5,4).Range.Text = "19+1=20"
$t.Cell(5,5).Range.Text = "64-28=36"

$t.Cell(6,1).Range.Text = "29+14=43"
$t.Cell(6,2).Range.Text = "97-50=47"
$t.Cell(6,3).Range.Text = "79+18=97"
$t.Cell(6,4).Range.Text = "54-26=28"
$t.Cell(6,5).Range.Text = "65-55=10"

$t.Cell(7,1).Range.Text = "1+81=82"
$t.Cell(7,2).Range.Text = "51+36=87"
$t.Cell(7,3).Range.Text = "67-48=19"
$t.Cell(7,4).Range.Text = "0+45=45"
$t.Cell(7,5).Range.Text = "68-56=12"

$t.Cell(8,1).Range.Text = "79-65=14"
$t.Cell(8,2).Range.Text = "15+4=19"
$t.Cell(8,3).Range.Text = "87-4=83"
$t.Cell(8,4).Range.Text = "98-27=71"
$t.Cell(8,5).Range.Text = "78-41=37"

$t.Cell(9,1).Range.Text = "99-40=59"
$t.Cell(9,2).Range.Text = "59-41=18"
$t.Cell(9,3).Range.Text = "17-11=6"
$t.Cell(9,4).Range.Text = "34-25=9"
$t.Cell(9,5).Range.Text = "63-4=59"

$t.Cell(10,1).Range.Text = "46+52=98"
$t.Cell(10,2).Range.Text = "79-59=20"
$t.Cell(10,3).Range.Text = "2+5=7"
$t.Cell(10,4).Range.Text = "78-37=41"
$t.Cell(10,5).Range.Text = "19-5=14"

$t.Cell(11,1).Range.Text = "6+51=57"
$t.Cell(11,2).Range.Text = "54+22=76"
$t.Cell(11,3).Range.Text = "25-15=10"
$t.Cell(11,4).Range.Text = "73-39=34"
$t.Cell(11,5).Range.Text = "67-7=60"

$t.Cell(12,1).Range.Text = "82+14=96"
$t.Cell(12,2).Range.Text = "34+11=45"
$t.Cell(12,3).Range.Text = "26+64=90"
$t.Cell(12,4).Range.Text = "46-38=8"
$t.Cell(12,5).Range.Text = "45+17=62"

$t.Cell(13,1).Range.Text = "87-11=76"
$t.Cell(13,2).Range.Text = "50-7=43"
$t.Cell(13,3).Range.Text = "44+52=96"
$t.Cell(13,4).Range.Text = "23+48=71"
$t.Cell(13,5).Range.Text = "98-38=60"

$t.Cell(14,1).Range.Text = "92-27=65"
$t.Cell(14,2).Range.Text = "6-1=5"
$t.Cell(14,3).Range.Text = "53+21=74"
$t.Cell(14,4).Range.Text = "2+39=41"
$t.Cell(14,5).Range.Text = "58+10=68"

$t.Cell(15,1).Range.Text = "30+23=53"
$t.Cell(15,2).Range.Text = "75-48=27"
$t.Cell(15,3).Range.Text = "71-35=36"
$t.Cell(15,4).Range.Text = "93-57=36"
$t.Cell(15,5).Range.Text = "41+34=75"

$t.Cell(16,1).Range.Text = "8+91=99"
$t.Cell(16,2).Range.Text = "38-6=32"
$t.Cell(16,3).Range.Text = "2+90=92"
$t.Cell(16,4).Range.Text = "38+23=61"
$t.Cell(16,5).Range.Text = "50+20=70"

$t.Cell(17,1).Range.Text = "4+39=43"
$t.Cell(17,2).Range.Text = "54-15=39"
$t.Cell(17,3).Range.Text = "22-11=11"
$t.Cell(17,4).Range.Text = "50-25=25"
$t.Cell(17,5).Range.Text = "89-30=59"

$t.Cell(18,1).Range.Text = "74-47=27"
$t.Cell(18,2).Range.Text = "24+18=42"
$t.Cell(18,3).Range.Text = "21+54=75"
$t.Cell(18,4).Range.Text = "44+38=82"
$t.Cell(18,5).Range.Text = "13+43=56"

$t.Cell(19,1).Range.Text = "50+6=56"
$t.Cell(19,2).Range.Text = "71+11=82"
$t.Cell(19,3).Range.Text = "43+5=48"
$t.Cell(19,4).Range.Text = "59-54=5"
$t.Cell(19,5).Range.Text = "34+35=69"

$t.Cell(20,1).Range.Text = "72-27=45"
$t.Cell(20,2).Range.Text = "68-7=61"
$t.Cell(20,3).Range.Text = "91-37=54"
$t.Cell(20,4).Range.Text = "7+25=32"
$t.Cell(20,5).Range.Text = "87-82=5"

